$d = $word.ActiveDocument

function Replace-Literal($find, $replace) {
    $rng = $d.Content
    $found = $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Text = $replace
    }
}

Replace-Literal "One of our researchers will call you to speak to you at a time that is convenient for you." "Een van ons navorsers sal jou bel op 'n tyd wat geskik is vir jou."

Replace-Literal "Why have I been invited to the interview?" "Hoekom was ek na die onderhoud genooi?"

Replace-Literal "The principal investigators of this study are Prof Cathy Ward and Cindee Bruyns and the Co-investigator is Carly Katzef all from the University of Cape Town." "Die hoofnavorsers van hierdie studie is Prof. Cathy Ward en Cindee Bruyns, en die Mede-navorsers is Carly Katzef, almal van die Universiteit van Kaapstad."

Replace-Literal "I am okay with the interview being recorded." "Ek is oukei daarmee dat die onderhoud opgeneem word."
